$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad / "Changed") holds a date serial number that was
# bumped by one day (46075 -> 46076, i.e. 2026-02-22 -> 2026-02-23) for
# every data row (rows 2 through 513).
$ws.Range("C2:C513").Value = 46076
